$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix spelling: "commited" -> "committed" everywhere it appears.
$ws.Range("B4").Value = "committed"
$ws.Range("B5").Value = "committed"
$ws.Range("B12").Value = "committed"

# 2. Changed dashboard layout: apply the same (Normal) cell style used by
#    the rest of the table to the last two data rows, matching rows above.
$ws.Range("A13:B14").Style = "Normal"

# 3. Move the active selection to B12.
$ws.Range("B12").Select()
